$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

for ($row = 2; $row -le 31; $row++) {
    $cell = $ws.Range("BF$row")
    $cell.NumberFormat = "@"
    $cell.Value = "2008-06-06"
    $cell.ClearFormats()
}
